$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.998.66"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "2.354.07"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.686"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.98"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.48%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +15.25%  "

$ws.Range("E10").Value = "  +1.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +17.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.08%  "

$ws.Range("E14").Value = "  +1.41%  "

$ws.Range("D15").Value = "2.705.69"
$ws.Range("E15").Value = "  -0.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.916"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.09%  "

$ws.Range("D18").Value = "2.351.59"
$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").Value = "43.918.21"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("E20").Value = "  +1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.62%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "

$ws.Range("E27").Value = "  +15.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.77%  "

$ws.Range("E29").Value = "  +1.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.137"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0758"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.66%  "

$ws.Range("E36").Value = "  +4.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.87%  "

$ws.Range("E38").Value = "  -2.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0282"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.98%  "

$ws.Range("E41").Value = "  +20.18%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.38%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.94%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.107"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.43%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.34%  "

$ws.Range("E48").Value = "  +2.71%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.20%  "
